$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("parameter_constraints")

# Update the constraint text: the total upper bound changes from 90 to 89.8
$ws.Range("A2").Value = "Na2SO4 + Na3PO4 + NaNO3 + MgSO4 + MgCl2 + K(CH3COOH) `n  + KH2PO4 + CaCl2 + ZnSO4 + Zn(CH3COOH) + MnSO4 + FeSO4 <= 89.8"

# Move the sheet's active selection from A2 to A3
$ws.Activate() | Out-Null
$ws.Range("A3").Select() | Out-Null
